$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 466.66666  # H2: 436.25 -> 466.66666
$ws.Cells.Item(2, 9).Value = 400  # I2: 398 -> 400
$ws.Cells.Item(2, 11).Value = 400  # K2: 398 -> 400
$ws.Cells.Item(2, 13).Value = -287  # M2: -285 -> -287
$ws.Cells.Item(29, 8).Value = 1158.4445  # H29: 1623.5 -> 1158.4445
$ws.Cells.Item(29, 10).Value = 1734.3334  # J29: 2753 -> 1734.3334
$ws.Cells.Item(29, 12).Value = 5203.0002  # L29: 8259 -> 5203.0002
$ws.Cells.Item(29, 14).Value = -5765.0002  # N29: -8821 -> -5765.0002
$ws.Cells.Item(38, 8).Value = 2260.923  # H38: 2080 -> 2260.923
$ws.Cells.Item(38, 9).Value = 233.16667  # I38: 223.375 -> 233.16667
$ws.Cells.Item(38, 10).Value = 3999  # J38: 3936.625 -> 3999
$ws.Cells.Item(38, 11).Value = 699.50001  # K38: 670.125 -> 699.50001
$ws.Cells.Item(38, 12).Value = 11997  # L38: 11809.875 -> 11997
$ws.Cells.Item(38, 13).Value = -327.50001  # M38: -298.125 -> -327.50001
$ws.Cells.Item(38, 14).Value = -12741  # N38: -12553.875 -> -12741
$ws.Cells.Item(62, 8).Value = 4699.3  # H62: 5944.1113 -> 4699.3
$ws.Cells.Item(62, 9).Value = 4699.3  # I62: 5944.1113 -> 4699.3
$ws.Cells.Item(62, 11).Value = 4699.3  # K62: 5944.1113 -> 4699.3
$ws.Cells.Item(62, 13).Value = -4075.3  # M62: -5320.1113 -> -4075.3
$ws.Cells.Item(65, 8).Value = 4699.3  # H65: 5944.1113 -> 4699.3
$ws.Cells.Item(65, 9).Value = 4699.3  # I65: 5944.1113 -> 4699.3
$ws.Cells.Item(65, 11).Value = 23496.5  # K65: 29720.5565 -> 23496.5
$ws.Cells.Item(65, 13).Value = -20376.5  # M65: -26600.5565 -> -20376.5
$ws.Cells.Item(76, 8).Value = 11473.286  # H76: 9931.299999999999 -> 11473.286
$ws.Cells.Item(76, 9).Value = 12549  # I76: 10032.667 -> 12549
$ws.Cells.Item(76, 10).Value = 10039  # J76: 9779.25 -> 10039
$ws.Cells.Item(76, 11).Value = 12549  # K76: 10032.667 -> 12549
$ws.Cells.Item(76, 12).Value = 10039  # L76: 9779.25 -> 10039
$ws.Cells.Item(76, 13).Value = -12234  # M76: -9717.666999999999 -> -12234
$ws.Cells.Item(76, 14).Value = -10669  # N76: -10409.25 -> -10669
$ws.Cells.Item(79, 8).Value = 11473.286  # H79: 9931.299999999999 -> 11473.286
$ws.Cells.Item(79, 9).Value = 12549  # I79: 10032.667 -> 12549
$ws.Cells.Item(79, 10).Value = 10039  # J79: 9779.25 -> 10039
$ws.Cells.Item(79, 11).Value = 12549  # K79: 10032.667 -> 12549
$ws.Cells.Item(79, 12).Value = 10039  # L79: 9779.25 -> 10039
$ws.Cells.Item(79, 13).Value = -11457  # M79: -8940.666999999999 -> -11457
$ws.Cells.Item(79, 14).Value = -12223  # N79: -11963.25 -> -12223
$ws.Cells.Item(113, 8).Value = 17166.834  # H113: 29417.334 -> 17166.834
$ws.Cells.Item(113, 9).Value = 19600.4  # I113: 29417.334 -> 19600.4
$ws.Cells.Item(113, 10).Value = 4999  # J113: 0 -> 4999
$ws.Cells.Item(113, 11).Value = 19600.4  # K113: 29417.334 -> 19600.4
$ws.Cells.Item(113, 12).Value = 4999  # L113: 0 -> 4999
$ws.Cells.Item(113, 13).Value = -16346.4  # M113: -26163.334 -> -16346.4
$ws.Cells.Item(113, 14).Value = -11507  # N113: None -> -11507
$ws.Cells.Item(116, 8).Value = 4882.1665  # H116: 4371.4546 -> 4882.1665
$ws.Cells.Item(116, 9).Value = 4799.3335  # I116: 4260.75 -> 4799.3335
$ws.Cells.Item(116, 10).Value = 4965  # J116: 4666.6665 -> 4965
$ws.Cells.Item(116, 11).Value = 4799.3335  # K116: 4260.75 -> 4799.3335
$ws.Cells.Item(116, 12).Value = 4965  # L116: 4666.6665 -> 4965
$ws.Cells.Item(116, 13).Value = -1357.3335  # M116: -818.75 -> -1357.3335
$ws.Cells.Item(116, 14).Value = -11849  # N116: -11550.6665 -> -11849
$ws.Cells.Item(118, 8).Value = 533.1667  # H118: 549.7143 -> 533.1667
$ws.Cells.Item(118, 9).Value = 424.75  # I118: 474.66666 -> 424.75
$ws.Cells.Item(118, 10).Value = 750  # J118: 1000 -> 750
$ws.Cells.Item(118, 11).Value = 1274.25  # K118: 1423.99998 -> 1274.25
$ws.Cells.Item(118, 12).Value = 2250  # L118: 3000 -> 2250
$ws.Cells.Item(118, 13).Value = 382.75  # M118: 233.0000199999999 -> 382.75
$ws.Cells.Item(118, 14).Value = -5564  # N118: -6314 -> -5564
$ws.Cells.Item(138, 8).Value = 3272.7083  # H138: 3328.7917 -> 3272.7083
$ws.Cells.Item(138, 9).Value = 2289.6428  # I138: 2412.3076 -> 2289.6428
$ws.Cells.Item(138, 10).Value = 3677.5  # J138: 3669.2 -> 3677.5
$ws.Cells.Item(138, 11).Value = 6868.928400000001  # K138: 7236.9228 -> 6868.928400000001
$ws.Cells.Item(138, 12).Value = 11032.5  # L138: 11007.6 -> 11032.5
$ws.Cells.Item(138, 13).Value = -1728.928400000001  # M138: -2096.9228 -> -1728.928400000001
$ws.Cells.Item(138, 14).Value = -21312.5  # N138: -21287.6 -> -21312.5
$ws.Cells.Item(141, 8).Value = 3877.2  # H141: 3380.611 -> 3877.2
$ws.Cells.Item(141, 9).Value = 3594.4167  # I141: 3055.0667 -> 3594.4167
$ws.Cells.Item(141, 11).Value = 10783.2501  # K141: 9165.2001 -> 10783.2501
$ws.Cells.Item(141, 13).Value = -5603.250100000001  # M141: -3985.2001 -> -5603.250100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3669.7385  # H32: 3910.41 -> 3669.7385
$ws.Cells.Item(32, 9).Value = 1733.6666  # I32: 1888.2195 -> 1733.6666
$ws.Cells.Item(32, 10).Value = 8025.9  # J32: 8055.9 -> 8025.9
$ws.Cells.Item(32, 11).Value = 1733.6666  # K32: 1888.2195 -> 1733.6666
$ws.Cells.Item(32, 12).Value = 8025.9  # L32: 8055.9 -> 8025.9
$ws.Cells.Item(32, 13).Value = -1446.6666  # M32: -1601.2195 -> -1446.6666
$ws.Cells.Item(32, 14).Value = -8599.9  # N32: -8629.9 -> -8599.9
$ws.Cells.Item(61, 8).Value = 871947.4399999999  # H61: 898339.75 -> 871947.4399999999
$ws.Cells.Item(61, 9).Value = 3571.7407  # I61: 3670.6538 -> 3571.7407
$ws.Cells.Item(61, 11).Value = 3571.7407  # K61: 3670.6538 -> 3571.7407
$ws.Cells.Item(61, 13).Value = -3359.7407  # M61: -3458.6538 -> -3359.7407
$ws.Cells.Item(110, 8).Value = 4993.1665  # H110: 6417.1113 -> 4993.1665
$ws.Cells.Item(110, 9).Value = 5356.273  # I110: 7094.375 -> 5356.273
$ws.Cells.Item(110, 11).Value = 5356.273  # K110: 7094.375 -> 5356.273
$ws.Cells.Item(110, 13).Value = -3311.273  # M110: -5049.375 -> -3311.273
$ws.Cells.Item(118, 8).Value = 224999.5  # H118: 250000 -> 224999.5
$ws.Cells.Item(118, 10).Value = 224999.5  # J118: 250000 -> 224999.5
$ws.Cells.Item(118, 12).Value = 224999.5  # L118: 250000 -> 224999.5
$ws.Cells.Item(118, 14).Value = -228313.5  # N118: -253314 -> -228313.5
$ws.Cells.Item(132, 8).Value = 1154817.9  # H132: 1115032.5 -> 1154817.9
$ws.Cells.Item(132, 9).Value = 2353.52  # I132: 2303 -> 2353.52
$ws.Cells.Item(132, 11).Value = 7060.559999999999  # K132: 6909 -> 7060.559999999999
$ws.Cells.Item(132, 13).Value = -4530.559999999999  # M132: -4379 -> -4530.559999999999
$ws.Cells.Item(136, 8).Value = 871947.4399999999  # H136: 898339.75 -> 871947.4399999999
$ws.Cells.Item(136, 9).Value = 3571.7407  # I136: 3670.6538 -> 3571.7407
$ws.Cells.Item(136, 11).Value = 10715.2221  # K136: 11011.9614 -> 10715.2221
$ws.Cells.Item(136, 13).Value = -8165.222099999999  # M136: -8461.9614 -> -8165.222099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 949.25  # H86: 0 -> 949.25
$ws.Cells.Item(86, 9).Value = 949.5  # I86: 0 -> 949.5
$ws.Cells.Item(86, 10).Value = 949  # J86: 0 -> 949
$ws.Cells.Item(86, 11).Value = 949.5  # K86: 0 -> 949.5
$ws.Cells.Item(86, 12).Value = 949  # L86: 0 -> 949
$ws.Cells.Item(86, 13).Value = 173.5  # M86: None -> 173.5
$ws.Cells.Item(86, 14).Value = -3195  # N86: None -> -3195
$ws.Cells.Item(89, 8).Value = 949.25  # H89: 0 -> 949.25
$ws.Cells.Item(89, 9).Value = 949.5  # I89: 0 -> 949.5
$ws.Cells.Item(89, 10).Value = 949  # J89: 0 -> 949
$ws.Cells.Item(89, 11).Value = 4747.5  # K89: 0 -> 4747.5
$ws.Cells.Item(89, 12).Value = 4745  # L89: 0 -> 4745
$ws.Cells.Item(89, 13).Value = 868.5  # M89: None -> 868.5
$ws.Cells.Item(89, 14).Value = -15977  # N89: None -> -15977
$ws.Cells.Item(94, 8).Value = 1299.6428  # H94: 1061.9231 -> 1299.6428
$ws.Cells.Item(94, 9).Value = 1550.7142  # I94: 1093 -> 1550.7142
$ws.Cells.Item(94, 10).Value = 1048.5714  # J94: 1042.5 -> 1048.5714
$ws.Cells.Item(94, 11).Value = 1550.7142  # K94: 1093 -> 1550.7142
$ws.Cells.Item(94, 12).Value = 1048.5714  # L94: 1042.5 -> 1048.5714
$ws.Cells.Item(94, 13).Value = -1099.7142  # M94: -642 -> -1099.7142
$ws.Cells.Item(94, 14).Value = -1950.5714  # N94: -1944.5 -> -1950.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(8, 8).Value = 2248.75  # H8: 2773.75 -> 2248.75
$ws.Cells.Item(8, 10).Value = 2965  # J8: 3665 -> 2965
$ws.Cells.Item(8, 12).Value = 2965  # L8: 3665 -> 2965
$ws.Cells.Item(8, 14).Value = -3245  # N8: -3945 -> -3245
$ws.Cells.Item(22, 8).Value = 1233.9445  # H22: 1306.8235 -> 1233.9445
$ws.Cells.Item(22, 9).Value = 310.22223  # I22: 336.14285 -> 310.22223
$ws.Cells.Item(22, 10).Value = 2157.6667  # J22: 1986.3 -> 2157.6667
$ws.Cells.Item(22, 11).Value = 310.22223  # K22: 336.14285 -> 310.22223
$ws.Cells.Item(22, 12).Value = 2157.6667  # L22: 1986.3 -> 2157.6667
$ws.Cells.Item(22, 13).Value = 39.77776999999998  # M22: 13.85714999999999 -> 39.77776999999998
$ws.Cells.Item(22, 14).Value = -2857.6667  # N22: -2686.3 -> -2857.6667
$ws.Cells.Item(31, 8).Value = 72867.3  # H31: 79873.45 -> 72867.3
$ws.Cells.Item(31, 9).Value = 156400.92  # I31: 184774.1 -> 156400.92
$ws.Cells.Item(31, 10).Value = 21156  # J31: 22178.1 -> 21156
$ws.Cells.Item(31, 11).Value = 156400.92  # K31: 184774.1 -> 156400.92
$ws.Cells.Item(31, 12).Value = 21156  # L31: 22178.1 -> 21156
$ws.Cells.Item(31, 13).Value = -156105.92  # M31: -184479.1 -> -156105.92
$ws.Cells.Item(31, 14).Value = -21746  # N31: -22768.1 -> -21746
$ws.Cells.Item(34, 8).Value = 72867.3  # H34: 79873.45 -> 72867.3
$ws.Cells.Item(34, 9).Value = 156400.92  # I34: 184774.1 -> 156400.92
$ws.Cells.Item(34, 10).Value = 21156  # J34: 22178.1 -> 21156
$ws.Cells.Item(34, 11).Value = 156400.92  # K34: 184774.1 -> 156400.92
$ws.Cells.Item(34, 12).Value = 21156  # L34: 22178.1 -> 21156
$ws.Cells.Item(34, 13).Value = -156198.92  # M34: -184572.1 -> -156198.92
$ws.Cells.Item(34, 14).Value = -21560  # N34: -22582.1 -> -21560
$ws.Cells.Item(58, 8).Value = 20921.045  # H58: 15124.032 -> 20921.045
$ws.Cells.Item(58, 9).Value = 14377.6  # I58: 8385.777 -> 14377.6
$ws.Cells.Item(58, 10).Value = 22845.588  # J58: 17880.592 -> 22845.588
$ws.Cells.Item(58, 11).Value = 14377.6  # K58: 8385.777 -> 14377.6
$ws.Cells.Item(58, 12).Value = 22845.588  # L58: 17880.592 -> 22845.588
$ws.Cells.Item(58, 13).Value = -14174.6  # M58: -8182.777 -> -14174.6
$ws.Cells.Item(58, 14).Value = -23251.588  # N58: -18286.592 -> -23251.588
$ws.Cells.Item(62, 8).Value = 4350  # H62: 3800 -> 4350
$ws.Cells.Item(62, 9).Value = 4300  # I62: 3600 -> 4300
$ws.Cells.Item(62, 10).Value = 4400  # J62: 3933.3333 -> 4400
$ws.Cells.Item(62, 11).Value = 4300  # K62: 3600 -> 4300
$ws.Cells.Item(62, 12).Value = 4400  # L62: 3933.3333 -> 4400
$ws.Cells.Item(62, 13).Value = -3676  # M62: -2976 -> -3676
$ws.Cells.Item(62, 14).Value = -5648  # N62: -5181.3333 -> -5648
$ws.Cells.Item(65, 8).Value = 4350  # H65: 3800 -> 4350
$ws.Cells.Item(65, 9).Value = 4300  # I65: 3600 -> 4300
$ws.Cells.Item(65, 10).Value = 4400  # J65: 3933.3333 -> 4400
$ws.Cells.Item(65, 11).Value = 21500  # K65: 18000 -> 21500
$ws.Cells.Item(65, 12).Value = 22000  # L65: 19666.6665 -> 22000
$ws.Cells.Item(65, 13).Value = -18380  # M65: -14880 -> -18380
$ws.Cells.Item(65, 14).Value = -28240  # N65: -25906.6665 -> -28240
$ws.Cells.Item(94, 8).Value = 18883.334  # H94: 19375 -> 18883.334
$ws.Cells.Item(94, 9).Value = 18427.25  # I94: 19164.75 -> 18427.25
$ws.Cells.Item(94, 11).Value = 18427.25  # K94: 19164.75 -> 18427.25
$ws.Cells.Item(94, 13).Value = -17976.25  # M94: -18713.75 -> -17976.25
$ws.Cells.Item(107, 8).Value = 500936.22  # H107: 612298.4399999999 -> 500936.22
$ws.Cells.Item(107, 9).Value = 688550.3  # I107: 918131.7 -> 688550.3
$ws.Cells.Item(107, 11).Value = 688550.3  # K107: 918131.7 -> 688550.3
$ws.Cells.Item(107, 13).Value = -686630.3  # M107: -916211.7 -> -686630.3
$ws.Cells.Item(131, 8).Value = 5000  # H131: 0 -> 5000
$ws.Cells.Item(131, 9).Value = 5000  # I131: 0 -> 5000
$ws.Cells.Item(131, 11).Value = 5000  # K131: 0 -> 5000
$ws.Cells.Item(131, 13).Value = 40  # M131: None -> 40
$ws.Cells.Item(132, 8).Value = 79196760  # H132: 85796330 -> 79196760
$ws.Cells.Item(132, 9).Value = 4613.4  # I132: 4902.6665 -> 4613.4
$ws.Cells.Item(132, 11).Value = 13840.2  # K132: 14707.9995 -> 13840.2
$ws.Cells.Item(132, 13).Value = -11310.2  # M132: -12177.9995 -> -11310.2
$ws.Cells.Item(136, 8).Value = 20921.045  # H136: 15124.032 -> 20921.045
$ws.Cells.Item(136, 9).Value = 14377.6  # I136: 8385.777 -> 14377.6
$ws.Cells.Item(136, 10).Value = 22845.588  # J136: 17880.592 -> 22845.588
$ws.Cells.Item(136, 11).Value = 43132.8  # K136: 25157.331 -> 43132.8
$ws.Cells.Item(136, 12).Value = 68536.764  # L136: 53641.776 -> 68536.764
$ws.Cells.Item(136, 13).Value = -40582.8  # M136: -22607.331 -> -40582.8
$ws.Cells.Item(136, 14).Value = -73636.764  # N136: -58741.776 -> -73636.764

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 3216468.5  # H4: 2924117.5 -> 3216468.5
$ws.Cells.Item(4, 9).Value = 3361492.8  # I4: 3460358 -> 3361492.8
$ws.Cells.Item(4, 10).Value = 2201299.8  # J4: 1100899.9 -> 2201299.8
$ws.Cells.Item(4, 11).Value = 10084478.4  # K4: 10381074 -> 10084478.4
$ws.Cells.Item(4, 12).Value = 6603899.399999999  # L4: 3302699.7 -> 6603899.399999999
$ws.Cells.Item(4, 13).Value = -10084366.4  # M4: -10380962 -> -10084366.4
$ws.Cells.Item(4, 14).Value = -6604123.399999999  # N4: -3302923.7 -> -6604123.399999999
$ws.Cells.Item(131, 8).Value = 1373.78  # H131: 1371.495 -> 1373.78
$ws.Cells.Item(131, 10).Value = 1489.3448  # J131: 1488.0581 -> 1489.3448
$ws.Cells.Item(131, 12).Value = 4468.0344  # L131: 4464.1743 -> 4468.0344
$ws.Cells.Item(131, 14).Value = -14548.0344  # N131: -14544.1743 -> -14548.0344
$ws.Cells.Item(132, 8).Value = 1596325.1  # H132: 1516553.9 -> 1596325.1
$ws.Cells.Item(132, 9).Value = 1383.25  # I132: 1346.0769 -> 1383.25
$ws.Cells.Item(132, 11).Value = 12449.25  # K132: 12114.6921 -> 12449.25
$ws.Cells.Item(132, 13).Value = -9919.25  # M132: -9584.6921 -> -9919.25
$ws.Cells.Item(137, 8).Value = 11448  # H137: 11142.8 -> 11448
$ws.Cells.Item(137, 10).Value = 20800  # J137: 22375 -> 20800
$ws.Cells.Item(137, 12).Value = 62400  # L137: 67125 -> 62400
$ws.Cells.Item(137, 14).Value = -72600  # N137: -77325 -> -72600

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 150  # H5: 200 -> 150
$ws.Cells.Item(5, 10).Value = 100  # J5: 0 -> 100
$ws.Cells.Item(5, 12).Value = 100  # L5: 0 -> 100
$ws.Cells.Item(5, 14).Value = -324  # N5: None -> -324
$ws.Cells.Item(70, 8).Value = 10687.031  # H70: 10870.968 -> 10687.031
$ws.Cells.Item(70, 9).Value = 9505.546  # I70: 9720.809999999999 -> 9505.546
$ws.Cells.Item(70, 11).Value = 9505.546  # K70: 9720.809999999999 -> 9505.546
$ws.Cells.Item(70, 13).Value = -9235.546  # M70: -9450.809999999999 -> -9235.546
$ws.Cells.Item(73, 8).Value = 10687.031  # H73: 10870.968 -> 10687.031
$ws.Cells.Item(73, 9).Value = 9505.546  # I73: 9720.809999999999 -> 9505.546
$ws.Cells.Item(73, 11).Value = 9505.546  # K73: 9720.809999999999 -> 9505.546
$ws.Cells.Item(73, 13).Value = -8569.546  # M73: -8784.809999999999 -> -8569.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(109, 8).Value = 10980  # H109: 12666.667 -> 10980
$ws.Cells.Item(109, 9).Value = 5000  # I109: 0 -> 5000
$ws.Cells.Item(109, 10).Value = 12475  # J109: 12666.667 -> 12475
$ws.Cells.Item(109, 11).Value = 5000  # K109: 0 -> 5000
$ws.Cells.Item(109, 12).Value = 12475  # L109: 12666.667 -> 12475
$ws.Cells.Item(109, 14).Value = -15249  # N109: -15440.667 -> -15249
$ws.Cells.Item(109, 13).Value = -3613  # M109: None -> -3613
$ws.Cells.Item(132, 8).Value = 1201666.4  # H132: 1305854.6 -> 1201666.4
$ws.Cells.Item(132, 9).Value = 4600.8  # I132: 5333.3335 -> 4600.8
$ws.Cells.Item(132, 11).Value = 13802.4  # K132: 16000.0005 -> 13802.4
$ws.Cells.Item(132, 13).Value = -11272.4  # M132: -13470.0005 -> -11272.4
$ws.Cells.Item(136, 8).Value = 892506.25  # H136: 1003794 -> 892506.25
$ws.Cells.Item(136, 9).Value = 14635.588  # I136: 17350.785 -> 14635.588
$ws.Cells.Item(136, 10).Value = 1677969.5  # J136: 1771027.6 -> 1677969.5
$ws.Cells.Item(136, 11).Value = 43906.764  # K136: 52052.355 -> 43906.764
$ws.Cells.Item(136, 12).Value = 5033908.5  # L136: 5313082.800000001 -> 5033908.5
$ws.Cells.Item(136, 13).Value = -41356.764  # M136: -49502.355 -> -41356.764
$ws.Cells.Item(136, 14).Value = -5039008.5  # N136: -5318182.800000001 -> -5039008.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 8000000  # H3: 2646545.5 -> 8000000
$ws.Cells.Item(3, 9).Value = 8000000  # I3: 5801400 -> 8000000
$ws.Cells.Item(3, 10).Value = 0  # J3: 17500 -> 0
$ws.Cells.Item(3, 11).Value = 8000000  # K3: 5801400 -> 8000000
$ws.Cells.Item(3, 12).Value = 0  # L3: 17500 -> 0
$ws.Cells.Item(3, 13).Value = -7999886  # M3: -5801286 -> -7999886
$ws.Cells.Item(3, 14).ClearContents()  # N3 was -17728
$ws.Cells.Item(41, 8).Value = 20247.334  # H41: 20376.334 -> 20247.334
$ws.Cells.Item(41, 9).Value = 19989  # I41: 0 -> 19989
$ws.Cells.Item(41, 10).Value = 20376.5  # J41: 20376.334 -> 20376.5
$ws.Cells.Item(41, 11).Value = 19989  # K41: 0 -> 19989
$ws.Cells.Item(41, 12).Value = 20376.5  # L41: 20376.334 -> 20376.5
$ws.Cells.Item(41, 14).Value = -21156.5  # N41: -21156.334 -> -21156.5
$ws.Cells.Item(41, 13).Value = -19599  # M41: None -> -19599
$ws.Cells.Item(105, 8).Value = 0  # H105: 20000 -> 0
$ws.Cells.Item(105, 9).Value = 0  # I105: 20000 -> 0
$ws.Cells.Item(105, 11).Value = 0  # K105: 20000 -> 0
$ws.Cells.Item(105, 13).ClearContents()  # M105 was -16506
$ws.Cells.Item(132, 8).Value = 2446484  # H132: 2795815 -> 2446484
$ws.Cells.Item(132, 9).Value = 5940.706  # I132: 6573.1333 -> 5940.706
$ws.Cells.Item(132, 10).Value = 8373517.5  # J132: 9768920 -> 8373517.5
$ws.Cells.Item(132, 11).Value = 17822.118  # K132: 19719.3999 -> 17822.118
$ws.Cells.Item(132, 12).Value = 25120552.5  # L132: 29306760 -> 25120552.5
$ws.Cells.Item(132, 13).Value = -15292.118  # M132: -17189.3999 -> -15292.118
$ws.Cells.Item(132, 14).Value = -25125612.5  # N132: -29311820 -> -25125612.5
$ws.Cells.Item(136, 8).Value = 446904.97  # H136: 487169.3 -> 446904.97
$ws.Cells.Item(136, 9).Value = 6783  # I136: 7600.3335 -> 6783
$ws.Cells.Item(136, 10).Value = 887026.9399999999  # J136: 819178.6 -> 887026.9399999999
$ws.Cells.Item(136, 11).Value = 20349  # K136: 22801.0005 -> 20349
$ws.Cells.Item(136, 12).Value = 2661080.82  # L136: 2457535.8 -> 2661080.82
$ws.Cells.Item(136, 13).Value = -17799  # M136: -20251.0005 -> -17799
$ws.Cells.Item(136, 14).Value = -2666180.82  # N136: -2462635.8 -> -2666180.82
